# Fixed some bugs in height
# Rows 2-23 (columns A-F) contain reel-weighting data that was reordered.
# Apply the corrected values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(301, 6, 45, 30, 60, 45)
    3  = @(601, 9, 60, 67, 60, 42)
    4  = @(901, 16, 15, 45, 60, 60)
    5  = @(801, 3, 67, 65, 52, 45)
    6  = @(501, 9, 52, 30, 75, 45)
    7  = @(201, 9, 30, 15, 45, 30)
    8  = @(1201, 2, 10, 10, 10, 10)
    9  = @(1202, 2, 10, 10, 10, 10)
    10 = @(1001, 18, 30, 75, 60, 72)
    11 = @(902, 1, 0, 0, 0, 0)
    12 = @(401, 9, 48, 67, 75, 45)
    13 = @(701, 3, 90, 45, 97, 15)
    14 = @(1203, 3, 15, 15, 15, 15)
    15 = @(101, 9, 30, 15, 60, 15)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(802, 0, 4, 5, 4, 0)
    18 = @(1101, 0, 15, 30, 30, 0)
    19 = @(1, 0, 2, 2, 2, 2)
    20 = @(2, 0, 2, 2, 2, 2)
    21 = @(3, 0, 3, 3, 3, 3)
    22 = @(602, 0, 0, 4, 0, 9)
    23 = @(402, 0, 0, 4, 0, 0)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
